# Scheduled runner update: refresh cached market-board pricing figures
# (currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfit columns)
# across the per-job leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 36520
$ws.Range("J93").Value = 36520
$ws.Range("L93").Value = 36520
$ws.Range("N93").Value = -41512

$ws.Range("H98").Value = 2910
$ws.Range("I98").Value = 2442.1428
$ws.Range("J98").Value = 4001.6667
$ws.Range("K98").Value = 2442.1428
$ws.Range("L98").Value = 4001.6667
$ws.Range("M98").Value = -944.1428000000001
$ws.Range("N98").Value = -6997.6667

$ws.Range("H111").Value = 73521.5
$ws.Range("J111").Value = 201700.2
$ws.Range("L111").Value = 605100.6000000001
$ws.Range("N111").Value = -611234.6000000001

$ws.Range("H122").Value = 2910
$ws.Range("I122").Value = 2442.1428
$ws.Range("J122").Value = 4001.6667
$ws.Range("K122").Value = 7326.428400000001
$ws.Range("L122").Value = 12005.0001
$ws.Range("M122").Value = -4876.428400000001
$ws.Range("N122").Value = -16905.0001

$ws.Range("H132").Value = 1536.125
$ws.Range("I132").Value = 1131.6786
$ws.Range("K132").Value = 3395.0358
$ws.Range("M132").Value = -865.0357999999997

$ws.Range("H138").Value = 2249.9583
$ws.Range("I138").Value = 849.7234
$ws.Range("J138").Value = 4882.4
$ws.Range("K138").Value = 2549.1702
$ws.Range("L138").Value = 14647.2
$ws.Range("M138").Value = 2590.8298
$ws.Range("N138").Value = -24927.2

$ws.Range("H141").Value = 925.35297
$ws.Range("I141").Value = 668.62067
$ws.Range("J141").Value = 2414.4
$ws.Range("K141").Value = 2005.86201
$ws.Range("L141").Value = 7243.200000000001
$ws.Range("M141").Value = 3174.13799
$ws.Range("N141").Value = -17603.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1207.95
$ws.Range("I32").Value = 1036.0804
$ws.Range("J32").Value = 2358.1538
$ws.Range("K32").Value = 1036.0804
$ws.Range("L32").Value = 2358.1538
$ws.Range("M32").Value = -749.0804000000001
$ws.Range("N32").Value = -2932.1538

$ws.Range("H61").Value = 1385.9556
$ws.Range("I61").Value = 1240.4117
$ws.Range("J61").Value = 1835.8182
$ws.Range("K61").Value = 1240.4117
$ws.Range("L61").Value = 1835.8182
$ws.Range("M61").Value = -1028.4117
$ws.Range("N61").Value = -2259.8182

$ws.Range("H136").Value = 1385.9556
$ws.Range("I136").Value = 1240.4117
$ws.Range("J136").Value = 1835.8182
$ws.Range("K136").Value = 3721.2351
$ws.Range("L136").Value = 5507.4546
$ws.Range("M136").Value = -1171.2351
$ws.Range("N136").Value = -10607.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 42155
$ws.Range("J76").Value = 42155
$ws.Range("L76").Value = 42155
$ws.Range("N76").Value = -42785

$ws.Range("H79").Value = 42155
$ws.Range("J79").Value = 42155
$ws.Range("L79").Value = 42155
$ws.Range("N79").Value = -44339

$ws.Range("H88").Value = 20760.334
$ws.Range("J88").Value = 20760.334
$ws.Range("L88").Value = 20760.334
$ws.Range("N88").Value = -21572.334

$ws.Range("H91").Value = 20760.334
$ws.Range("J91").Value = 20760.334
$ws.Range("L91").Value = 20760.334
$ws.Range("N91").Value = -23568.334

$ws.Range("H134").Value = 2145.5
$ws.Range("I134").Value = 1862.0476
$ws.Range("J134").Value = 2428.9524
$ws.Range("K134").Value = 5586.142800000001
$ws.Range("L134").Value = 7286.8572
$ws.Range("M134").Value = -3051.142800000001
$ws.Range("N134").Value = -12356.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 40780
$ws.Range("J52").Value = 40780
$ws.Range("L52").Value = 40780
$ws.Range("N52").Value = -41368

$ws.Range("H74").Value = 26699.75
$ws.Range("J74").Value = 26699.75
$ws.Range("L74").Value = 26699.75
$ws.Range("N74").Value = -28447.75

$ws.Range("H77").Value = 26699.75
$ws.Range("J77").Value = 26699.75
$ws.Range("L77").Value = 80099.25
$ws.Range("N77").Value = -88835.25

$ws.Range("H92").Value = 45000
$ws.Range("J92").Value = 45000
$ws.Range("L92").Value = 45000
$ws.Range("N92").Value = -49992

$ws.Range("H132").Value = 1958.898
$ws.Range("I132").Value = 1590.1389
$ws.Range("J132").Value = 2980.077
$ws.Range("K132").Value = 4770.4167
$ws.Range("L132").Value = 8940.231
$ws.Range("M132").Value = -2240.4167
$ws.Range("N132").Value = -14000.231

$ws.Range("H134").Value = 2189.4614
$ws.Range("I134").Value = 2556.2424
$ws.Range("J134").Value = 1552.421
$ws.Range("K134").Value = 7668.7272
$ws.Range("L134").Value = 4657.263
$ws.Range("M134").Value = -5133.7272
$ws.Range("N134").Value = -9727.262999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3333447.2
$ws.Range("I12").Value = 7692390
$ws.Range("J12").Value = 138.11765
$ws.Range("K12").Value = 23077170
$ws.Range("L12").Value = 414.35295
$ws.Range("M12").Value = -23076997
$ws.Range("N12").Value = -760.35295

$ws.Range("H113").Value = 167182.08
$ws.Range("I113").Value = 518.74286
$ws.Range("J113").Value = 400510.75
$ws.Range("K113").Value = 1556.22858
$ws.Range("L113").Value = 1201532.25
$ws.Range("M113").Value = 613.77142
$ws.Range("N113").Value = -1205872.25

$ws.Range("H132").Value = 1612.375
$ws.Range("I132").Value = 640
$ws.Range("K132").Value = 5760
$ws.Range("M132").Value = -3230

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 4031
$ws.Range("J47").Value = 4031
$ws.Range("L47").Value = 4031
$ws.Range("N47").Value = -5167

$ws.Range("H126").Value = 5846.1724
$ws.Range("I126").Value = 11842.4
$ws.Range("J126").Value = 2690.2632
$ws.Range("K126").Value = 35527.2
$ws.Range("L126").Value = 8070.7896
$ws.Range("M126").Value = -33057.2
$ws.Range("N126").Value = -13010.7896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1704407.5
$ws.Range("I122").Value = 2236178.5
$ws.Range("J122").Value = 2739.8
$ws.Range("K122").Value = 6708535.5
$ws.Range("L122").Value = 8219.400000000001
$ws.Range("M122").Value = -6706085.5
$ws.Range("N122").Value = -13119.4

$ws.Range("H132").Value = 10277755
$ws.Range("I132").Value = 11617826
$ws.Range("J132").Value = 3880
$ws.Range("K132").Value = 34853478
$ws.Range("L132").Value = 11640
$ws.Range("M132").Value = -34850948
$ws.Range("N132").Value = -16700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7464933
$ws.Range("I136").Value = 2466.182
$ws.Range("J136").Value = 21740956
$ws.Range("K136").Value = 7398.545999999999
$ws.Range("L136").Value = 65222868
$ws.Range("M136").Value = -4848.545999999999
$ws.Range("N136").Value = -65227968
